$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.621.54'
$ws.Range("E2").Value = '  +3.86%  '
$ws.Range("D3").Value = '2.634.41'
$ws.Range("E3").Value = '  +2.48%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.176'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.79%  '
$ws.Range("D10").Value = '2.632.10'
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("E12").Value = '  +3.41%  '
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("E14").Value = '  +4.28%  '
$ws.Range("D15").Value = '3.108.01'
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("D16").Value = '72.478.80'
$ws.Range("E16").Value = '  +3.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.81'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").Value = '2.620.12'
$ws.Range("E18").Value = '  +2.06%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '385.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.29%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("E23").Value = '  +15.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.70%  '
$ws.Range("E25").Value = '  +2.70%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.27%  '
$ws.Range("D28").Value = '2.730.20'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("D30").Value = '0.0₃0962'
$ws.Range("E30").Value = '  +4.17%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '522.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.24%  '
$ws.Range("E33").Value = '  +3.44%  '
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.43'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("E38").Value = '  +3.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.98%  '
$ws.Range("E40").Value = '  -5.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.30%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  +4.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.334'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("E49").Value = '  +4.11%  '
$ws.Range("E50").Value = '  +4.84%  '
$ws.Range("D51").Value = '0.0₆0266'
$ws.Range("E51").Value = '  +2.16%  '
